$d = $word.ActiveDocument
$tbl = $d.Tables.Item(2)
$row = $tbl.Rows.Item(7)

$r1 = $row.Cells.Item(3).Range
$r1.Text = "Aadhar"
$r1.Font.Name = "Lora"

$r2 = $row.Cells.Item(4).Range
$r2.Text = "{aadhar}"
$r2.Font.Name = "Lora"
